$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.938.38"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "3.266.28"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "185.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "580.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "3.264.10"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("E10").Value = "  -4.05%  "
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").Value = "3.835.79"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.138"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.20%  "
$ws.Range("D16").Value = "67.954.22"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("E17").Value = "  -2.77%  "
$ws.Range("D18").Value = "3.288.91"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "397.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.99%  "
$ws.Range("E22").Value = "  -2.59%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("E26").Value = "  -4.24%  "
$ws.Range("E27").Value = "  -3.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.84%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("E32").Value = "  -6.82%  "
$ws.Range("E33").Value = "  -4.11%  "
$ws.Range("E34").Value = "  -5.91%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -5.33%  "
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.809"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.08%  "
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.71%  "
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").Value = "2.645.97"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("E46").Value = "  -8.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "333.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.82%  "
$ws.Range("E49").Value = "  -3.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.101"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.70%  "
